# Dodano losowanie wybranej przez admina liczby pytan (narzucone w kodzie).
# Adds a new "test" question bank (questions 7-10) to the Pytania sheet so
# the quiz has a larger pool to randomly draw the admin-chosen number of
# questions from.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 32
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "test"
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = "Która odp jest dobra?"
$ws.Range("E33").Value = "TA"
$ws.Range("E32").Value = "Nie TA"
$ws.Range("F32").Value = "F"
$ws.Range("G32").Value = "F"

# Row 33
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = "test"
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = "Która odp jest dobra?"
$ws.Range("F33").Value = "T"
$ws.Range("G33").Value = "F"

# Row 34
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "test"
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = "Która odp jest dobra?"
$ws.Range("E34").Value = "nie tu"
$ws.Range("F34").Value = "F"
$ws.Range("G34").Value = "F"

# Row 35
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "test"
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = "Która odp jest dobra?"
$ws.Range("E35").Value = "nie tam"
$ws.Range("F35").Value = "F"
$ws.Range("G35").Value = "F"

# Row 36
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "test"
$ws.Range("C36").Value = 2
$ws.Range("D36").Value = "Która odp jest dobra?"
$ws.Range("E36").Value = "nie nigdzie"
$ws.Range("F36").Value = "F"
$ws.Range("G36").Value = "F"

# Row 37
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "test"
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = "Która odp jest dobra?"
$ws.Range("E37").Value = "smutek"
$ws.Range("F37").Value = "F"
$ws.Range("G37").Value = "F"

# Row 38
$ws.Range("A38").Value = 8
$ws.Range("B38").Value = "test"
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = "Jaki był najlepszy utwor Jap poprzedniego sezonu?"
$ws.Range("E38").Value = "Silent solitude"
$ws.Range("F38").Value = "T"
$ws.Range("G38").Value = "F"

# Row 39
$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "test"
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = "Jaki był najlepszy utwor Jap poprzedniego sezonu?"
$ws.Range("E39").Value = "Innocent note"
$ws.Range("F39").Value = "F"
$ws.Range("G39").Value = "F"

# Row 40
$ws.Range("A40").Value = 8
$ws.Range("B40").Value = "test"
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = "Jaki był najlepszy utwor Jap poprzedniego sezonu?"
$ws.Range("E40").Value = "Marionette "
$ws.Range("F40").Value = "F"
$ws.Range("G40").Value = "F"

# Row 41
$ws.Range("A41").Value = 8
$ws.Range("B41").Value = "test"
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = "Jaki był najlepszy utwor Jap poprzedniego sezonu?"
$ws.Range("E41").Value = "Blue Bird "
$ws.Range("F41").Value = "F"
$ws.Range("G41").Value = "F"

# Row 42
$ws.Range("A42").Value = 8
$ws.Range("B42").Value = "test"
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = "Jaki był najlepszy utwor Jap poprzedniego sezonu?"
$ws.Range("E42").Value = "D Techno Life "
$ws.Range("F42").Value = "F"
$ws.Range("G42").Value = "F"

# Row 43
$ws.Range("A43").Value = 8
$ws.Range("B43").Value = "test"
$ws.Range("C43").Value = 2
$ws.Range("D43").Value = "Jaki był najlepszy utwor Jap poprzedniego sezonu?"
$ws.Range("E43").Value = "Adamas"
$ws.Range("F43").Value = "F"
$ws.Range("G43").Value = "F"

# Row 44
$ws.Range("A44").Value = 9
$ws.Range("B44").Value = "test"
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = "Sao to bajka?"
$ws.Range("E44").Value = "Super"
$ws.Range("F44").Value = "F"
$ws.Range("G44").Value = "F"

# Row 45
$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "test"
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = "Sao to bajka?"
$ws.Range("E45").Value = "Bardzo dobra"
$ws.Range("F45").Value = "F"
$ws.Range("G45").Value = "F"

# Row 46
$ws.Range("A46").Value = 9
$ws.Range("B46").Value = "test"
$ws.Range("C46").Value = 3
$ws.Range("D46").Value = "Sao to bajka?"
$ws.Range("E46").NumberFormat = "d-mmm"
$ws.Range("E46").Value = "piec na dziesiec"
$ws.Range("F46").Value = "F"
$ws.Range("G46").Value = "F"

# Row 47
$ws.Range("A47").Value = 9
$ws.Range("B47").Value = "test"
$ws.Range("C47").Value = 3
$ws.Range("D47").Value = "Sao to bajka?"
$ws.Range("E47").Value = "Powinna się skonczyc "
$ws.Range("F47").Value = "F"
$ws.Range("G47").Value = "F"

# Row 48
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "test"
$ws.Range("C48").Value = 3
$ws.Range("D48").Value = "Sao to bajka?"
$ws.Range("E48").Value = "Dlaczego ono istnieje"
$ws.Range("F48").Value = "F"
$ws.Range("G48").Value = "F"

# Row 49
$ws.Range("A49").Value = 9
$ws.Range("B49").Value = "test"
$ws.Range("C49").Value = 3
$ws.Range("D49").Value = "Sao to bajka?"
$ws.Range("E49").Value = "Cos się popsulo i nie było go slychac"
$ws.Range("F49").Value = "T"
$ws.Range("G49").Value = "F"

# Row 50
$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "test"
$ws.Range("C50").Value = 3
$ws.Range("D50").Value = "Sao to bajka?"
$ws.Range("E50").Value = "Chyba ok."
$ws.Range("F50").Value = "F"
$ws.Range("G50").Value = "F"

# Row 51
$ws.Range("A51").Value = 9
$ws.Range("B51").Value = "test"
$ws.Range("C51").Value = 3
$ws.Range("D51").Value = "Sao to bajka?"
$ws.Range("E51").Value = "Trudne pytanka"
$ws.Range("F51").Value = "F"
$ws.Range("G51").Value = "F"

# Row 52
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "test"
$ws.Range("C52").Value = 3
$ws.Range("D52").Value = "Losowanie to ciezka sprawa"
$ws.Range("E52").Value = "Bardzo"
$ws.Range("F52").Value = "F"
$ws.Range("G52").Value = "F"

# Row 53
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "test"
$ws.Range("C53").Value = 3
$ws.Range("D53").Value = "Losowanie to ciezka sprawa"
$ws.Range("E53").Value = "Bayes "
$ws.Range("F53").Value = "T"
$ws.Range("G53").Value = "F"

# Row 54
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "test"
$ws.Range("C54").Value = 3
$ws.Range("D54").Value = "Losowanie to ciezka sprawa"
$ws.Range("E54").Value = "Wszystko losowe "
$ws.Range("F54").Value = "F"
$ws.Range("G54").Value = "F"

# Row 55
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "test"
$ws.Range("C55").Value = 3
$ws.Range("D55").Value = "Losowanie to ciezka sprawa"
$ws.Range("E55").Value = "Nic losowe, wszystko dane"
$ws.Range("F55").Value = "F"
$ws.Range("G55").Value = "F"

# Row 56
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "test"
$ws.Range("C56").Value = 3
$ws.Range("D56").Value = "Losowanie to ciezka sprawa"
$ws.Range("E56").Value = "Rownania rozniczkowe"
$ws.Range("F56").Value = "F"
$ws.Range("G56").Value = "F"

# Row 57
$ws.Range("A57").Value = 10
$ws.Range("B57").Value = "test"
$ws.Range("C57").Value = 3
$ws.Range("D57").Value = "Losowanie to ciezka sprawa"
$ws.Range("E57").Value = "KHUN TACKER"
$ws.Range("F57").Value = "F"
$ws.Range("G57").Value = "F"

# Row 58
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "test"
$ws.Range("C58").Value = 3
$ws.Range("D58").Value = "Losowanie to ciezka sprawa"
$ws.Range("E58").Value = "Model dynamiczny arrowa hurtowicza "
$ws.Range("F58").Value = "F"
$ws.Range("G58").Value = "F"

# Row 59
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "test"
$ws.Range("C59").Value = 3
$ws.Range("D59").Value = "Losowanie to ciezka sprawa"
$ws.Range("E59").Value = "kek, herbata się zrobila "
$ws.Range("F59").Value = "F"
$ws.Range("G59").Value = "F"

$ws.Range("E61").Select()
